$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: add hanging indent (w:ind w:left="720" w:hanging="720")
# ---------------------------------------------------------------------
$pTitle = $d.Paragraphs(1)
$pTitle.LeftIndent = 36        # 720 twips = 36pt
$pTitle.FirstLineIndent = -36  # hanging indent

# ---------------------------------------------------------------------
# 2. Fix typo "refatoriserat" -> "refaktoriserat"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("refatoriserat", $false, $false, $false, $false, $false, $true, 1, $false, "refaktoriserat", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Replace the tab after "...nästkommande vecka." with a space, merging
#    it into the following sentence ("Vi har även implementerat...").
# ---------------------------------------------------------------------
$tab = [char]9
$search  = "nästkommande vecka." + $tab + "Vi har även"
$replace = "nästkommande vecka. Vi har även"
$d.Content.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Build the skeleton of new paragraphs after paragraph 7 ("Zombies...").
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8.Style = "Heading1"

$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs(9)
$p9.Style = "Normal"

$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs(10)
$p10.Style = "Normal"

$p10.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs(11)
$p11.Style = "Normal"

$p11.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs(12)
$p12.Style = "Normal"

$p12.Range.InsertParagraphAfter()
$p13 = $d.Paragraphs(13)
$p13.Style = "Normal"

$p13.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs(14)
$p14.Style = "Normal"

$p14.Range.InsertParagraphAfter()
$p15 = $d.Paragraphs(15)
$p15.Style = "Normal"

# ---------------------------------------------------------------------
# 5. Paragraph 8: "V. 6." heading
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8.Range.InsertAfter("V. 6.")

# ---------------------------------------------------------------------
# 6. Paragraph 9 stays empty.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 7. Paragraph 10: "Då vi har arbetat..." with a line break and an
#    italic "Inventory" run.
# ---------------------------------------------------------------------
$vt = [char]11   # manual line break char, serializes to <w:br/>

$part1 = "Då vi har arbetat med mycket olika features det senaste men inte riktigt hunnit klart med allt för att kunna merge:a in det i develop har det därför hänt en hel del den här veckan. Vi har under veckans gång vävt ihop alla trådar som har legat ute i periferin, och har nu en konkret applikation som börjar likna vårt slutliga mål."
$part2 = $vt + "   Under veckans gång har vi gjort klart allt som har med hälsa att göra, man kan nu alltså ta skada och bli skadad av zombies. Människor (då även spelaren själv) blir nu infekterad då en zombie attackerar denne. "
$part3 = "Vapenförrådet ("

$p10 = $d.Paragraphs(10)
$p10.Range.InsertAfter($part1 + $part2 + $part3)

$p10 = $d.Paragraphs(10)
$italicStart = $p10.Range.End - 1
$p10.Range.InsertAfter("Inventory")
$p10 = $d.Paragraphs(10)
$italicEnd = $p10.Range.End - 1
$rItalic = $d.Range($italicStart, $italicEnd)
$rItalic.Italic = 1

$part4 = ") är nu även färdigställt med animationer som vi känner oss nöjda med. Man kan nu flytta vapen från ett större"
$part5 = " "
$part6 = "inventory,"
$part7 = " som dyker upp då man pausar spelet"
$part8 = ", till ett mer lättåtkomligt inventory (går enkelt att komma åt under spelets gång). "

$p10 = $d.Paragraphs(10)
$p10.Range.InsertAfter($part4 + $part5 + $part6 + $part7 + $part8)

# ---------------------------------------------------------------------
# 8. Paragraph 11: start-menu / buy zones paragraph.
# ---------------------------------------------------------------------
$p11Text1 = "   En enklare startmeny startas nu även då man drar igång spelet (så länge inte debug-mode är aktiverat), med tillhörande bilder. "
$p11Text2 = "Vi har även gjort klart logiken för hur det ska gå till då man som spelare går in i en byggnad, dvs att hustaket ska försvinna. Det finns nu även specifika ”buy zones” på banan där man kan gå in för att köpa nya vapen som man kan använda sig av under den desperata kampen mot de hjärnfrossande zombies som vandrar runt i staden. "

$p11 = $d.Paragraphs(11)
$p11.Range.InsertAfter($p11Text1 + $p11Text2)

# ---------------------------------------------------------------------
# 9. Paragraph 12: just three spaces.
# ---------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$p12.Range.InsertAfter("   ")

# ---------------------------------------------------------------------
# 10. Paragraph 13: final summary paragraph + relocated _GoBack bookmark.
# ---------------------------------------------------------------------
$p13Text = "Utöver ovanstående är nu även vapen med tillhörande bilder fixade, man tar alltså olika mycket hälsa från fienderna då man skjuter på dem med olika vapen. Ljud är även implementerat i spelet för att skapa stämning, något vi kommer att fortsätta att jobba med nästkommande vecka. Sedan har det även legat fokus på att förbättra kodstrukturen och prestandan för applikationen."

$p13 = $d.Paragraphs(13)
$p13.Range.InsertAfter($p13Text)

# Move the _GoBack bookmark from paragraph 7 to the end of paragraph 13.
$p13 = $d.Paragraphs(13)
$bmPoint = $d.Range($p13.Range.End - 1, $p13.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ---------------------------------------------------------------------
# 11. Paragraph 14 stays empty; paragraph 15 gets "both" justification.
# ---------------------------------------------------------------------
$p15 = $d.Paragraphs(15)
$p15.Alignment = 3   # wdAlignParagraphJustify -> <w:jc w:val="both"/>

# ---------------------------------------------------------------------
# 12. Recompute pagination so any required lastRenderedPageBreak markers
#     are (re)generated.
# ---------------------------------------------------------------------
try {
    $d.Repaginate()
} catch {
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
